$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1641221374045801
$ws.Range("C2").Value = 0.6221374045801527
$ws.Range("J2").Value = 0.007633587786259542
$ws.Range("P2").Value = 0.1145038167938931
$ws.Range("S2").Value = 0.0916030534351145
$ws.Range("B3").Value = 0.01176470588235294
$ws.Range("C3").Value = 0.04705882352941176
$ws.Range("J3").Value = 0.01176470588235294
$ws.Range("P3").Value = 0.6941176470588235
$ws.Range("S3").Value = 0.2352941176470588
$ws.Range("J4").Value = 0.03278688524590164
$ws.Range("O4").Value = 0.01639344262295082
$ws.Range("P4").Value = 0.7049180327868853
$ws.Range("S4").Value = 0.2459016393442623
$ws.Range("B6").Value = 0.0660377358490566
$ws.Range("D6").Value = 0.02358490566037736
$ws.Range("F6").Value = 0.08490566037735849
$ws.Range("J6").Value = 0.2735849056603774
$ws.Range("O6").Value = 0.01415094339622642
$ws.Range("Q6").Value = 0.1981132075471698
$ws.Range("R6").Value = 0.0330188679245283
$ws.Range("S6").Value = 0.3066037735849056
$ws.Range("B7").Value = 0.1344086021505376
$ws.Range("D7").Value = 0.04301075268817205
$ws.Range("E7").Value = 0.005376344086021506
$ws.Range("F7").Value = 0.04838709677419355
$ws.Range("J7").Value = 0.06451612903225806
$ws.Range("O7").Value = 0.02688172043010753
$ws.Range("Q7").Value = 0.1827956989247312
$ws.Range("R7").Value = 0.06451612903225806
$ws.Range("S7").Value = 0.4301075268817204
$ws.Range("B8").Value = 0.09264305177111716
$ws.Range("D8").Value = 0.01907356948228883
$ws.Range("E8").Value = 0.002724795640326975
$ws.Range("F8").Value = 0.07356948228882834
$ws.Range("J8").Value = 0.1035422343324251
$ws.Range("O8").Value = 0.03542234332425068
$ws.Range("Q8").Value = 0.2397820163487738
$ws.Range("R8").Value = 0.06267029972752043
$ws.Range("S8").Value = 0.3705722070844686
$ws.Range("B9").Value = 0.1016042780748663
$ws.Range("D9").Value = 0.0160427807486631
$ws.Range("F9").Value = 0.0374331550802139
$ws.Range("J9").Value = 0.08021390374331551
$ws.Range("O9").Value = 0.0213903743315508
$ws.Range("Q9").Value = 0.2406417112299465
$ws.Range("R9").Value = 0.05882352941176471
$ws.Range("S9").Value = 0.4438502673796791
$ws.Range("B10").Value = 0.1200369344413666
$ws.Range("D10").Value = 0.03785780240073869
$ws.Range("E10").Value = 0.0009233610341643582
$ws.Range("F10").Value = 0.07017543859649122
$ws.Range("J10").Value = 0.1089566020313943
$ws.Range("O10").Value = 0.02862419205909511
$ws.Range("Q10").Value = 0.221606648199446
$ws.Range("R10").Value = 0.04986149584487535
$ws.Range("S10").Value = 0.3619575253924284
$ws.Range("G11").Value = 0.1517241379310345
$ws.Range("J11").Value = 0.07931034482758621
$ws.Range("K11").Value = 0.1862068965517241
$ws.Range("L11").Value = 0.5724137931034483
$ws.Range("S11").Value = 0.0103448275862069
$ws.Range("G12").Value = 0.7251461988304093
$ws.Range("J12").Value = 0.2163742690058479
$ws.Range("K12").Value = 0.005847953216374269
$ws.Range("L12").Value = 0.01169590643274854
$ws.Range("S12").Value = 0.04093567251461988
$ws.Range("G13").Value = 0.5652173913043478
$ws.Range("J13").Value = 0.3695652173913043
$ws.Range("S13").Value = 0.06521739130434782
$ws.Range("F15").Value = 0.02202643171806168
$ws.Range("H15").Value = 0.09691629955947137
$ws.Range("I15").Value = 0.07929515418502203
$ws.Range("J15").Value = 0.3436123348017621
$ws.Range("K15").Value = 0.07929515418502203
$ws.Range("M15").Value = 0.01762114537444934
$ws.Range("O15").Value = 0.1013215859030837
$ws.Range("S15").Value = 0.2599118942731278
$ws.Range("F16").Value = 0.0267379679144385
$ws.Range("H16").Value = 0.1550802139037433
$ws.Range("I16").Value = 0.09625668449197861
$ws.Range("J16").Value = 0.4064171122994653
$ws.Range("K16").Value = 0.1390374331550802
$ws.Range("M16").Value = 0.0267379679144385
$ws.Range("O16").Value = 0.053475935828877
$ws.Range("S16").Value = 0.09625668449197861
$ws.Range("F17").Value = 0.02921348314606742
$ws.Range("H17").Value = 0.1730337078651685
$ws.Range("I17").Value = 0.09887640449438202
$ws.Range("J17").Value = 0.3887640449438202
$ws.Range("K17").Value = 0.1303370786516854
$ws.Range("M17").Value = 0.02471910112359551
$ws.Range("O17").Value = 0.0651685393258427
$ws.Range("S17").Value = 0.0898876404494382
$ws.Range("H18").Value = 0.09345794392523364
$ws.Range("I18").Value = 0.102803738317757
$ws.Range("J18").Value = 0.4392523364485981
$ws.Range("K18").Value = 0.102803738317757
$ws.Range("M18").Value = 0.02803738317757009
$ws.Range("O18").Value = 0.08411214953271028
$ws.Range("S18").Value = 0.1495327102803738
$ws.Range("F19").Value = 0.02033271719038817
$ws.Range("H19").Value = 0.2134935304990758
$ws.Range("I19").Value = 0.08687615526802218
$ws.Range("J19").Value = 0.3696857670979667
$ws.Range("K19").Value = 0.1090573012939002
$ws.Range("M19").Value = 0.02125693160813309
$ws.Range("N19").Value = 0.0009242144177449168
$ws.Range("O19").Value = 0.07116451016635859
$ws.Range("S19").Value = 0.1072088724584103
